$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so Excel does not
# coerce them into numbers (they are text values like "42.84").
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D5').Value = '233.51'
$ws.Range('D7').Value = '0.4691'
$ws.Range('D8').Value = '42.84'
$ws.Range('D9').Value = '0.2824'
$ws.Range('D10').Value = '0.06445'
$ws.Range('D11').Value = '20.93'
$ws.Range('D12').Value = '0.07734'
$ws.Range('D14').Value = '93.28'
$ws.Range('D15').Value = '0.6765'
$ws.Range('D16').Value = '5.033'
$ws.Range('D17').Value = '266.09'
$ws.Range('D19').Value = '13.28'
$ws.Range('D20').Value = '0.000007551'
$ws.Range('D24').Value = '5.131'
$ws.Range('D25').Value = '6.084'
$ws.Range('D26').Value = '9.272'
$ws.Range('D27').Value = '165.25'
$ws.Range('D29').Value = '1.876'
$ws.Range('D30').Value = '1.363'
$ws.Range('D31').Value = '0.09814'
$ws.Range('D32').Value = '1.450'
$ws.Range('D33').Value = '4.185'
$ws.Range('D34').Value = '3.962'
$ws.Range('D35').Value = '0.04632'
$ws.Range('D37').Value = '0.6849'
$ws.Range('D38').Value = '2.714'
$ws.Range('D39').Value = '0.01826'
$ws.Range('D40').Value = '2.712'
$ws.Range('D41').Value = '6.238'
$ws.Range('D42').Value = '70.41'
$ws.Range('D44').Value = '0.8309'
$ws.Range('D45').Value = '101.69'
$ws.Range('D46').Value = '1.859'
$ws.Range('D47').Value = '0.4028'
$ws.Range('D48').Value = '9.111'
$ws.Range('D49').Value = '6.911'
$ws.Range('D50').Value = '914.60'
$ws.Range('D51').Value = '33.94'

# Reset style of Price column back to default (remove the quote-prefix/
# text-format artifact) so only the value text differs, matching source.
$ws.Range("D2:D51").Style = "Normal"

# Remaining cells (Volume % text, Coin names, Links) are already text and
# can be set directly without the NumberFormat workaround.
$ws.Range('D2').Value = '30.105.05'
$ws.Range('E2').Value = '  -0.61%  '
$ws.Range('D3').Value = '1.856.51'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('E5').Value = '  -0.80%  '
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('E9').Value = '  -1.71%  '
$ws.Range('E10').Value = '  -2.16%  '
$ws.Range('E11').Value = '  -4.20%  '
$ws.Range('D13').Value = '1.871.36'
$ws.Range('E13').Value = '  -0.04%  '
$ws.Range('E14').Value = '  -4.07%  '
$ws.Range('E15').Value = '  -1.37%  '
$ws.Range('E16').Value = '  -1.86%  '
$ws.Range('E17').Value = '  -1.18%  '
$ws.Range('D18').Value = '30.082.20'
$ws.Range('E18').Value = '  -0.64%  '
$ws.Range('E19').Value = '  -5.54%  '
$ws.Range('E20').Value = '  -1.58%  '
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('D22').Value = '2.113.48'
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('E24').Value = '  -2.78%  '
$ws.Range('E25').Value = '  -2.21%  '
$ws.Range('E26').Value = '  -1.68%  '
$ws.Range('E27').Value = '  -1.78%  '
$ws.Range('E28').Value = '  -2.47%  '
$ws.Range('E29').Value = '  -3.83%  '
$ws.Range('E30').Value = '  -0.43%  '
$ws.Range('E32').Value = '  -0.97%  '
$ws.Range('E33').Value = '  -4.49%  '
$ws.Range('E34').Value = '  -2.91%  '
$ws.Range('E35').Value = '  -1.67%  '
$ws.Range('E36').Value = '  -2.20%  '
$ws.Range('E37').Value = '  -2.34%  '
$ws.Range('E38').Value = '  +0.09%  '
$ws.Range('E39').Value = '  -2.67%  '
$ws.Range('E40').Value = '  +3.35%  '
$ws.Range('E41').Value = '  -1.10%  '
$ws.Range('E42').Value = '  -2.61%  '
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('E44').Value = '  -1.52%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('E45').Value = '  -1.33%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('E46').Value = '  -4.92%  '
$ws.Range('E47').Value = '  -3.32%  '
$ws.Range('E48').Value = '  -0.69%  '
$ws.Range('E49').Value = '  -2.20%  '
$ws.Range('E50').Value = '  -1.17%  '
$ws.Range('E51').Value = '  -1.73%  '
